# Reorder the comma-separated email lists in column G ("Recorded By")
# for the rows affected by the sync, without altering any other content,
# formatting, or formulas on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$updates = @{
    2  = "Veronia.rafat@med.asu.edu.eg, servinaz@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, gehanadel@med.asu.edu.eg, System"
    3  = "Veronia.rafat@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, System, hend_mahmoud@med.asu.edu.eg"
    4  = "servinaz@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, gehanadel@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg"
    5  = "asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"
    7  = "lamiaa.ossama@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, Kerelos.zareef@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg, Fatmaelhady@med.asu.edu.eg"
    9  = "Shimaa.ashraf@med.asu.edu.eg, Safa.hany@med.asu.edu.eg"
    12 = "Madeha.Saeed@med.asu.edu.eg, dina.adel@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg"
    15 = "Rania.a.youssef@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"
    28 = "Aya_hamed@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg"
}

foreach ($row in $updates.Keys) {
    $ws.Range("G$row").Value = $updates[$row]
}
